# Re-export of the speaker-variant playlist: rows are re-ordered (no
# preferred-speaker-variant ordering, no Levenshtein-distance grouping
# applied anymore), and the "is_prefered" column is no longer populated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data occupies rows 2..50 (row 1 is the header).
$firstRow = 2
$lastRow = 50
$rowCount = $lastRow - $firstRow + 1

# --- 1. Snapshot the current (id / speaker_variant) pairs for every data row
#        before anything gets overwritten. ---------------------------------
$ids = @()
$variants = @()
for ($i = $firstRow; $i -le $lastRow; $i++) {
    $ids += , ($ws.Cells.Item($i, 2).Value())
    $variants += , ($ws.Cells.Item($i, 3).Value())
}

# --- 2. New row order: for each destination row (index 0 == row 2) this
#        gives the 1-based offset (within the snapshot above) of the row
#        whose (id, speaker_variant) pair should land there. -----------------
$order = @(8,20,18,25,50,7,23,6,22,30,27,15,37,5,14,31,21,29,3,46,34,36,45,4,33,11,41,42,35,39,26,28,19,10,44,40,47,12,49,9,16,32,13,17,43,38,48,24,2)

for ($i = 0; $i -lt $rowCount; $i++) {
    $destRow = $firstRow + $i
    $srcOffset = $order[$i] - $firstRow
    $ws.Cells.Item($destRow, 2).Value = $ids[$srcOffset]
    $ws.Cells.Item($destRow, 3).Value = $variants[$srcOffset]

    # The "is_prefered" column is no longer exported - clear it everywhere.
    $ws.Cells.Item($destRow, 4).Value = ""
}

# --- 3. A handful of ids are re-slugified from their (unchanged) display
#        text now that the Levenshtein-based matching is gone. -------------
$ws.Range("B6").Value = "#e-vrouw"
$ws.Range("B38").Value = "#zuylesteyn"
$ws.Range("B40").Value = "#natuur--en-staat--beschouwers"
$ws.Range("B48").Value = "#coenraat"
